$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 164, shifting the existing rows 164..263
# down to 166..265 (Excel copies formatting, incl. the date style, down
# automatically).
$ws.Rows("164:165").Insert()

# --- New row 164 ---
$ws.Range("A164").Value = 10
$ws.Range("B164").Value = "Vega Modelo de Temuco"
$ws.Range("C164").Value = "La Araucanía"
$ws.Range("D164").Value = (Get-Date -Year 2022 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E164").Value = 9
$ws.Range("F164").Value = 100114013
$ws.Range("G164").Value = "Zanahoria"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 150
$ws.Range("K164").Value = 8000
$ws.Range("L164").Value = 8000
$ws.Range("M164").Value = 8000
$ws.Range("N164").Value = "`$/saco 20 kilos"
$ws.Range("O164").Value = "Región del Maule"
$ws.Range("P164").Value = 400
$ws.Range("Q164").Value = 20
$ws.Range("R164").Value = "Hortaliza"

# --- New row 165 ---
$ws.Range("A165").Value = 10
$ws.Range("B165").Value = "Vega Modelo de Temuco"
$ws.Range("C165").Value = "La Araucanía"
$ws.Range("D165").Value = (Get-Date -Year 2022 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E165").Value = 9
$ws.Range("F165").Value = 100114013
$ws.Range("G165").Value = "Zanahoria"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 100
$ws.Range("K165").Value = 8000
$ws.Range("L165").Value = 8000
$ws.Range("M165").Value = 8000
$ws.Range("N165").Value = "`$/saco 25 kilos"
$ws.Range("O165").Value = "Región de La Araucanía"
$ws.Range("P165").Value = 320
$ws.Range("Q165").Value = 25
$ws.Range("R165").Value = "Hortaliza"
